# Weekly update: insert a new price record for "Ají" (Feria Lagunitas de
# Puerto Montt) as the new row 127, pushing the existing rows 127:149 down
# to 128:150 (dimension grows from A1:R149 to A1:R150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 127 (and everything below it) down by one row.
$ws.Rows(127).Insert()

# Populate the newly inserted row 127 with this week's record.
$ws.Range("A127").Value = 4
$ws.Range("B127").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C127").Value = "Los Lagos"
$ws.Range("D127").Value = 44476
$ws.Range("E127").Value = 10
$ws.Range("F127").Value = 100112021
$ws.Range("G127").Value = "Ají"
$ws.Range("H127").Value = "Inferno"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 60
$ws.Range("K127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("M127").Value = 50000
$ws.Range("N127").Value = "`$/caja 12 kilos"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 4167
$ws.Range("Q127").Value = 12
$ws.Range("R127").Value = "Hortaliza"
